$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.261.82"
$ws.Range("E2").Value = "  -2.82%  "
$ws.Range("D3").Value = "1.556.98"
$ws.Range("E3").Value = "  -4.29%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "206.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.19%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.243"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.98%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.0609"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "17.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.0781"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").Value = "1.772.08"
$ws.Range("E12").Value = "  -4.37%  "
$ws.Range("D13").Value = "1.551.97"
$ws.Range("E13").Value = "  -4.79%  "
$ws.Range("E14").Value = "  -4.14%  "
$ws.Range("E15").Value = "  -3.92%  "
$ws.Range("D16").Value = "25.259.24"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "59.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.94%  "
$ws.Range("E18").Value = "  -4.37%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("E20").Value = "  -3.33%  "
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "5.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.34%  "
$ws.Range("E24").Value = "  -3.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "140.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("E27").Value = "  -4.27%  "
$ws.Range("E28").Value = "  -2.36%  "
$ws.Range("E29").Value = "  -4.51%  "
$ws.Range("E30").Value = "  -7.04%  "
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "3.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("E33").Value = "  -4.61%  "
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("E35").Value = "  -3.81%  "
$ws.Range("D36").Value = "1.085.50"
$ws.Range("E36").Value = "  -3.36%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "2.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.16%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.0149"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.493"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.766"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.800"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "92.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.65%  "
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").Value = "1.686.74"
$ws.Range("E45").Value = "  -4.33%  "
$ws.Range("E46").Value = "  -2.50%  "
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "52.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.64%  "
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("E51").Value = "  -2.25%  "
